# "Added year to archetypes (only one construction per function for now)"
#
# For every building-use row (rows 2-19) on both the ARCHITECTURE and HVAC
# sheets, fill in the previously-empty "year_start" (B) and "year_end" (C)
# columns with the value 0 (stored as text, matching shared-string "0").
#
# Also mirrors the author's change of active sheet/selection: the HVAC tab
# becomes the active one (instead of ARCHITECTURE), with its selection
# anchored on B19:C19, while ARCHITECTURE's own selection moves to C19.

$wb = $excel.ActiveWorkbook

$wsArchitecture = $wb.Worksheets.Item("ARCHITECTURE")
$wsHvac = $wb.Worksheets.Item("HVAC")

# Fill year_start (B) / year_end (C) for rows 2-19 with 0 on both sheets.
$wsArchitecture.Range("B2:C19").Value = "0"
$wsHvac.Range("B2:C19").Value = "0"

# Update the lingering selection on ARCHITECTURE before switching away from it.
[void]$wsArchitecture.Range("C19").Select()

# HVAC becomes the active/selected sheet, with its own new selection.
[void]$wsHvac.Activate()
[void]$wsHvac.Range("B19:C19").Select()
